# Rename column headers from "_old"/"_new" suffixes to the respective
# format-version suffixes ("_FV2410" / "_FV2504"), turn the header row +
# data range into a proper Excel Table ("Table1"), and freeze the header
# row (top row) so it stays visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header cells -------------------------------------------
$oldSuffix = "_old"
$newSuffix = "_FV2410"
$oldSuffix2 = "_new"
$newSuffix2 = "_FV2504"

for ($col = 1; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $text = $cell.Value()
    if ($text -eq "diff") {
        continue
    }
    if ($text -like "*$oldSuffix") {
        $base = $text.Substring(0, $text.Length - $oldSuffix.Length)
        $cell.Value = "$base$newSuffix"
    } elseif ($text -like "*$oldSuffix2") {
        $base = $text.Substring(0, $text.Length - $oldSuffix2.Length)
        $cell.Value = "$base$newSuffix2"
    }
}

# --- 2. Convert the used range into an Excel Table ------------------------
$usedRange = $ws.Range("A1:U78")
$table = $ws.ListObjects.Add(1, $usedRange, 0, 1)
$table.Name = "Table1"

# --- 3. Freeze the header (top) row ----------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
